$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new data values in column C
$ws.Range("C2").Value = 96
$ws.Range("C4").Value = 36

# Update the selected cell (active cell) to H14
$ws.Range("H14").Select()
